# Auto-generated edit script: applies numeric value updates to the
# "Twintania_Profits" leve-profit workbook, sheet by sheet, row by row,
# matching the target OOXML diff (scheduled-runner data refresh).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 12150.934
$ws.Range("I6").Value = 12150.934
$ws.Range("K6").Value = 36452.802
$ws.Range("M6").Value = -36340.802
$ws.Range("H41").Value = 793.8570999999999
$ws.Range("I41").Value = 1077
$ws.Range("J41").Value = 333.75
$ws.Range("K41").Value = 1077
$ws.Range("L41").Value = 333.75
$ws.Range("M41").Value = -637
$ws.Range("N41").Value = -1213.75
$ws.Range("H112").Value = 3448.7
$ws.Range("I112").Value = 1414.8334
$ws.Range("K112").Value = 4244.5002
$ws.Range("M112").Value = -3136.5002
$ws.Range("H132").Value = 2859087.8
$ws.Range("I132").Value = 3391409
$ws.Range("J132").Value = 3910
$ws.Range("K132").Value = 10174227
$ws.Range("L132").Value = 11730
$ws.Range("M132").Value = -10171697
$ws.Range("N132").Value = -16790
$ws.Range("H137").Value = 14122.223
$ws.Range("I137").Value = 7213.875
$ws.Range("J137").Value = 27938.916
$ws.Range("K137").Value = 21641.625
$ws.Range("L137").Value = 83816.74800000001
$ws.Range("M137").Value = -19091.625
$ws.Range("N137").Value = -88916.74800000001

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9631.538
$ws.Range("I45").Value = 12880.223
$ws.Range("J45").Value = 2322
$ws.Range("K45").Value = 12880.223
$ws.Range("L45").Value = 2322
$ws.Range("M45").Value = -12503.223
$ws.Range("N45").Value = -3076
$ws.Range("H61").Value = 3574.372
$ws.Range("I61").Value = 2329.3784
$ws.Range("K61").Value = 2329.3784
$ws.Range("M61").Value = -2117.3784
$ws.Range("H63").Value = 1325.5
$ws.Range("I63").Value = 639.44446
$ws.Range("J63").Value = 7500
$ws.Range("K63").Value = 639.44446
$ws.Range("L63").Value = 7500
$ws.Range("M63").Value = 46.55553999999995
$ws.Range("N63").Value = -8872
$ws.Range("H66").Value = 1325.5
$ws.Range("I66").Value = 639.44446
$ws.Range("J66").Value = 7500
$ws.Range("K66").Value = 3197.2223
$ws.Range("L66").Value = 37500
$ws.Range("M66").Value = 234.7776999999996
$ws.Range("N66").Value = -44364
$ws.Range("H122").Value = 2148
$ws.Range("I122").Value = 2148
$ws.Range("K122").Value = 6444
$ws.Range("M122").Value = -3994
$ws.Range("H136").Value = 3574.372
$ws.Range("I136").Value = 2329.3784
$ws.Range("K136").Value = 6988.135200000001
$ws.Range("M136").Value = -4438.135200000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2927.8
$ws.Range("I99").Value = 2999.3635
$ws.Range("K99").Value = 2999.3635
$ws.Range("M99").Value = -1501.3635

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31253000
$ws.Range("J31").Value = 4500
$ws.Range("L31").Value = 4500
$ws.Range("N31").Value = -5090
$ws.Range("H34").Value = 31253000
$ws.Range("J34").Value = 4500
$ws.Range("L34").Value = 4500
$ws.Range("N34").Value = -4904
$ws.Range("H58").Value = 21595.285
$ws.Range("I58").Value = 22168.777
$ws.Range("K58").Value = 22168.777
$ws.Range("M58").Value = -21965.777
$ws.Range("H99").Value = 38991.383
$ws.Range("I99").Value = 89779.664
$ws.Range("K99").Value = 89779.664
$ws.Range("M99").Value = -88281.664
$ws.Range("H126").Value = 38991.383
$ws.Range("I126").Value = 89779.664
$ws.Range("K126").Value = 269338.992
$ws.Range("M126").Value = -266868.992
$ws.Range("H136").Value = 21595.285
$ws.Range("I136").Value = 22168.777
$ws.Range("K136").Value = 66506.33099999999
$ws.Range("M136").Value = -63956.33099999999

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 62500136
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H68").Value = 1330.3846
$ws.Range("J68").Value = 1472.125
$ws.Range("L68").Value = 4416.375
$ws.Range("N68").Value = -6038.375
$ws.Range("H71").Value = 1330.3846
$ws.Range("J71").Value = 1472.125
$ws.Range("L71").Value = 13249.125
$ws.Range("N71").Value = -21361.125
$ws.Range("H92").Value = 88.666664
$ws.Range("J92").Value = 83.5
$ws.Range("L92").Value = 250.5
$ws.Range("N92").Value = -2746.5
$ws.Range("H109").Value = 4778.5
$ws.Range("I109").Value = 4778.5
$ws.Range("K109").Value = 14335.5
$ws.Range("M109").Value = -13295.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2175.6667
$ws.Range("I122").Value = 1810.8
$ws.Range("K122").Value = 5432.4
$ws.Range("M122").Value = -2982.4
$ws.Range("H123").Value = 30326
$ws.Range("J123").Value = 30326
$ws.Range("L123").Value = 30326
$ws.Range("N123").Value = -35226

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1213.2667
$ws.Range("I22").Value = 1477.091
$ws.Range("K22").Value = 1477.091
$ws.Range("M22").Value = -1182.091
$ws.Range("H27").Value = 1213.2667
$ws.Range("I27").Value = 1477.091
$ws.Range("K27").Value = 1477.091
$ws.Range("M27").Value = -1370.091
$ws.Range("H43").Value = 99500
$ws.Range("I43").Value = 100000
$ws.Range("J43").Value = 99000
$ws.Range("K43").Value = 100000
$ws.Range("L43").Value = 99000
$ws.Range("M43").Value = -99807
$ws.Range("N43").Value = -99386
$ws.Range("H46").Value = 1282.2059
$ws.Range("I46").Value = 931.3125
$ws.Range("J46").Value = 1594.1111
$ws.Range("K46").Value = 931.3125
$ws.Range("L46").Value = 1594.1111
$ws.Range("M46").Value = -743.3125
$ws.Range("N46").Value = -1970.1111
$ws.Range("H55").Value = 287.8421
$ws.Range("I55").Value = 285.6
$ws.Range("J55").Value = 296.25
$ws.Range("K55").Value = 285.6
$ws.Range("L55").Value = 296.25
$ws.Range("M55").Value = -112.6
$ws.Range("N55").Value = -642.25
$ws.Range("H132").Value = 7154.1
$ws.Range("I132").Value = 5075.9165
$ws.Range("J132").Value = 15466.833
$ws.Range("K132").Value = 15227.7495
$ws.Range("L132").Value = 46400.499
$ws.Range("M132").Value = -12697.7495
$ws.Range("N132").Value = -51460.499

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 16409.59
$ws.Range("I132").Value = 14287.929
$ws.Range("J132").Value = 21810.182
$ws.Range("K132").Value = 42863.787
$ws.Range("L132").Value = 65430.546
$ws.Range("M132").Value = -40333.787
$ws.Range("N132").Value = -70490.546
$ws.Range("H136").Value = 2626
$ws.Range("I136").Value = 2283.4333
$ws.Range("J136").Value = 3482.4167
$ws.Range("K136").Value = 6850.2999
$ws.Range("L136").Value = 10447.2501
$ws.Range("N136").Value = -15547.2501
